# repull data, push all data, mean calculation
# Update the "dSF" column (F) for a set of rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    7  = 3
    13 = -5
    14 = 7
    21 = -6
    25 = -7
    31 = -3
    35 = 0
    36 = -2
    40 = 0
    43 = 3
    44 = -8
    48 = -4
    51 = 3
    54 = -5
    57 = -3
    60 = -2
    62 = 1
    63 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
